$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy header style from C1 to D1 and E1 (bold, centered, bordered)
$ws.Range("C1").Copy($ws.Range("D1:E1"))

# Update header texts
$ws.Range("C1").Value = "Frecuencia del primer armonico"
$ws.Range("D1").Value = "Frecuencia del segundo armonico"
$ws.Range("E1").Value = "Frecuencia tercer armonico"

# Update data rows 2-48 for columns C, D, E
$ws.Range("C2").Value = 301.2429929320006
$ws.Range("D2").Value = 149.1591518401169
$ws.Range("E2").Value = 456.2515232756523
$ws.Range("C3").Value = 182.6992473450873
$ws.Range("D3").Value = 366.6357356428498
$ws.Range("E3").Value = 549.7473966388288
$ws.Range("C4").Value = 393.69267236937
$ws.Range("D4").Value = 262.5991961249097
$ws.Range("E4").Value = 660.0020612181797
$ws.Range("C5").Value = 330.2127659574467
$ws.Range("D5").Value = 165.1063829787236
$ws.Range("E5").Value = 495.3191489361707
$ws.Range("C6").Value = 149.5527180977138
$ws.Range("D6").Value = 297.576267298723
$ws.Range("E6").Value = 450.4931569691871
$ws.Range("C7").Value = 470.2338485929449
$ws.Range("D7").Value = 234.3242172017435
$ws.Range("E7").Value = 704.2409829567978
$ws.Range("C8").Value = 375.4609453570229
$ws.Range("D8").Value = 187.0600067046594
$ws.Range("E8").Value = 748.2400268186393
$ws.Range("C9").Value = 301.3663644489843
$ws.Range("D9").Value = 150.6831822244922
$ws.Range("E9").Value = 458.1790320521004
$ws.Range("C10").Value = 206.2364927446743
$ws.Range("D10").Value = 407.9448389420604
$ws.Range("E10").Value = 621.1793763507258
$ws.Range("C11").Value = 296.1519727228442
$ws.Range("D11").Value = 593.6028576067542
$ws.Range("E11").Value = 370.1899659035557
$ws.Range("C12").Value = 415.6487873959991
$ws.Range("D12").Value = 217.3836077181804
$ws.Range("E12").Value = 623.8272260577091
$ws.Range("C13").Value = 356.2610229276897
$ws.Range("D13").Value = 179.0123456790125
$ws.Range("E13").Value = 704.5855379188715
$ws.Range("C14").Value = 185.5742296918761
$ws.Range("D14").Value = 374.6498599439774
$ws.Range("E14").Value = 556.0224089635849
$ws.Range("C15").Value = 330.5275637225841
$ws.Range("D15").Value = 165.5008891523412
$ws.Range("E15").Value = 496.5026674570245
$ws.Range("C16").Value = 206.276612521684
$ws.Range("D16").Value = 411.9224097145561
$ws.Range("E16").Value = 618.409294012511
$ws.Range("C17").Value = 358.4229390681003
$ws.Range("D17").Value = 179.8088410991641
$ws.Range("E17").Value = 545.4002389486259
$ws.Range("C18").Value = 369.0219068188831
$ws.Range("D18").Value = 185.5394425588811
$ws.Range("E18").Value = 732.6956700606811
$ws.Range("C19").Value = 345.7663072299806
$ws.Range("D19").Value = 487.1840197984798
$ws.Range("E19").Value = 171.8225207707264
$ws.Range("C20").Value = 348.5707672723065
$ws.Range("D20").Value = 174.5168383288974
$ws.Range("E20").Value = 523.0876056012039
$ws.Range("C21").Value = 326.7197682838519
$ws.Range("D21").Value = 490.6589427950767
$ws.Range("E21").Value = 163.3598841419257
$ws.Range("C22").Value = 341.4094701920176
$ws.Range("D22").Value = 173.3185885191515
$ws.Range("E22").Value = 511.9131396400926
$ws.Range("C23").Value = 302.6679524268729
$ws.Range("D23").Value = 453.359048537448
$ws.Range("E23").Value = 150.9482481517198
$ws.Range("C24").Value = 317.264957264957
$ws.Range("D24").Value = 476.5811965811963
$ws.Range("E24").Value = 158.6324786324785
$ws.Range("C25").Value = 294.3905525094897
$ws.Range("D25").Value = 571.9105862505276
$ws.Range("E25").Value = 858.7094053142137
$ws.Range("C26").Value = 264.0867150407594
$ws.Range("D26").Value = 529.2484099256471
$ws.Range("E26").Value = 797.2767177282094
$ws.Range("C27").Value = 294.1803453421444
$ws.Range("D27").Value = 591.7714772969516
$ws.Range("E27").Value = 885.951822639096
$ws.Range("C28").Value = 256.8161829375549
$ws.Range("D28").Value = 769.666764389719
$ws.Range("E28").Value = 514.0232580865822
$ws.Range("C29").Value = 541.017347838871
$ws.Range("D29").Value = 269.3325492502208
$ws.Range("E29").Value = 785.6512790355778
$ws.Range("C30").Value = 489.5241824946152
$ws.Range("D30").Value = 245.153710593303
$ws.Range("E30").Value = 738.594086547876
$ws.Range("C31").Value = 269.2180798416366
$ws.Range("D31").Value = 539.0960079181787
$ws.Range("E31").Value = 806.9943912900035
$ws.Range("C32").Value = 422.8471281765314
$ws.Range("D32").Value = 211.2642396239944
$ws.Range("E32").Value = 851.7485859953795
$ws.Range("C33").Value = 233.5913157101213
$ws.Range("D33").Value = 466.7805809629108
$ws.Range("E33").Value = 699.9698462157003
$ws.Range("C34").Value = 243.8858695652179
$ws.Range("D34").Value = 489.8097826086955
$ws.Range("E34").Value = 728.260869565218
$ws.Range("C35").Value = 287.2027180067953
$ws.Range("D35").Value = 578.9354473386184
$ws.Range("E35").Value = 862.514156285391
$ws.Range("C36").Value = 290.6326849988818
$ws.Range("D36").Value = 578.5826067516209
$ws.Range("E36").Value = 862.0612564274534
$ws.Range("C37").Value = 274.2481384779039
$ws.Range("D37").Value = 827.386132869161
$ws.Range("E37").Value = 549.6567063146695
$ws.Range("C38").Value = 193.1420573827845
$ws.Range("D38").Value = 375.6473058082574
$ws.Range("E38").Value = 564.8705388383478
$ws.Range("C39").Value = 293.5740839086566
$ws.Range("D39").Value = 586.7233138608608
$ws.Range("E39").Value = 880.7222517259697
$ws.Range("C40").Value = 558.6516473278562
$ws.Range("D40").Value = 279.4952147031418
$ws.Range("E40").Value = 836.4529516388584
$ws.Range("C41").Value = 311.970454051705
$ws.Range("D41").Value = 629.1548989789262
$ws.Range("E41").Value = 935.0423636758633
$ws.Range("C42").Value = 584.160813756358
$ws.Range("D42").Value = 876.0797610397999
$ws.Range("E42").Value = 289.6585129571322
$ws.Range("C43").Value = 170.5041917581766
$ws.Range("D43").Value = 343.3699374188218
$ws.Range("E43").Value = 514.3464399574923
$ws.Range("C44").Value = 432.0015553611347
$ws.Range("D44").Value = 288.1306503353744
$ws.Range("E44").Value = 144.2597453096141
$ws.Range("C45").Value = 393.9111525318422
$ws.Range("D45").Value = 198.8195091643365
$ws.Range("E45").Value = 801.4911463187318
$ws.Range("C46").Value = 193.1206614148659
$ws.Range("D46").Value = 385.9293346852819
$ws.Range("E46").Value = 768.4267997816078
$ws.Range("C47").Value = 435.4213089048444
$ws.Range("D47").Value = 218.8660559544442
$ws.Range("E47").Value = 346.9505653214492
$ws.Range("C48").Value = 250.1745635910229
$ws.Range("D48").Value = 499.1521197007487
$ws.Range("E48").Value = 748.927680798005
